$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.711.22'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +8.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.384.38'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.74%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '418.92'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +6.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '115.92'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +8.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.597'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +6.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.649'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +5.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.78'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.108'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +10.24%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.925.82'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.55'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +5.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.10'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +6.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.474.74'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +7.21%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.494.75'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +8.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.86'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000116'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +10.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.30'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '306.82'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.97'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.24'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.11'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +4.64%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.36%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.79'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +8.30%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.02'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.180'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +7.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.116'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +7.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.58'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +22.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.59'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +6.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '39.95'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +7.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0514'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +6.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.62'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.15'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.44'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '137.75'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.01%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.294'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.90%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.05'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.12'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.30'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +11.92%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.80'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.171.62'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.41'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.02'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.40%  '
